# Generate Report for Handback
# Updates the handback status timestamps / status text that are refreshed
# each time the report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# "Latest HO Xliff Generate Date" for a14d245e...md / dbfedc10...md rows
# (shared between the Overview sheet and the de-de sheet's
# "Correspond Handoff Datetime" column for the same rows)
$wsOverview.Range("G3").Value = "2016-08-27 18:23:22"
$wsOverview.Range("G5").Value = "2016-08-27 18:23:22"
$wsDeDe.Range("H3").Value = "2016-08-27 18:23:22"
$wsDeDe.Range("H5").Value = "2016-08-27 18:23:22"

# Status column: "ht" -> "mt" for the a14d245e.../dbfedc10...md rows,
# on both the zh-cn and de-de sheets
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn "Correspond Handoff Datetime" for the same rows
$wsZhCn.Range("H3").Value = "2016-08-27 18:23:17"
$wsZhCn.Range("H5").Value = "2016-08-27 18:23:17"

# zh-cn "Correspond Handback DateTime" for the same rows
$wsZhCn.Range("K3").Value = "2016-08-27 18:23:41"
$wsZhCn.Range("K5").Value = "2016-08-27 18:23:41"

# de-de "Correspond Handback DateTime" for the same rows
$wsDeDe.Range("K3").Value = "2016-08-27 18:23:48"
$wsDeDe.Range("K5").Value = "2016-08-27 18:23:48"
